# Error Calculations and Plots
# This script reproduces the edits made to the missing_data.xlsx worksheet:
#  - Two data rows ("RM 232" and "SC 92") were removed from the table,
#    shifting all subsequent rows up.
#  - A number of cells in columns D and E (imputed/missing-value columns)
#    were updated: some previously-blank cells received newly imputed
#    values, and some previously-filled cells were cleared back to blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the two rows that were dropped from the table ---
# Row 26 = "RM 232" (first, since it appears earlier in the sheet)
$ws.Rows("26").Delete()
# After the above deletion, the row that used to be "SC 92" (row 28)
# has shifted up to row 27.
$ws.Rows("27").Delete()

# --- Apply the individual D/E column value updates (post row-deletion row numbers) ---

# Row 2 "RM 2": D=missing already set; E blank -> -7.2
$ws.Range("E2").Value = -7.2

# Row 6 "RM 21": E -5.7 -> blank
$ws.Range("E6").ClearContents()

# Row 12 "RM 81": E blank -> -5.3
$ws.Range("E12").Value = -5.3

# Row 14 "RM 90": E -5.4 -> blank
$ws.Range("E14").ClearContents()

# Row 20 "RM 134": E blank -> -7.2
$ws.Range("E20").Value = -7.2

# Row 21 "RM 135": E blank -> -8.699999999999999
$ws.Range("E21").Value = -8.699999999999999

# Row 23 "RM 140": E -7 -> blank
$ws.Range("E23").ClearContents()

# Row 24 "RM 142a": E -8.1 -> blank
$ws.Range("E24").ClearContents()

# Row 26 "SC 5": D blank -> -13.8
$ws.Range("D26").Value = -13.8

# Row 27 "SC 101": D -14.6 -> blank
$ws.Range("D27").ClearContents()

# Row 28 "SC 105": D -13.7 -> blank
$ws.Range("D28").ClearContents()

# Row 29 "SC 119": D blank -> -13
$ws.Range("D29").Value = -13

# Row 30 "SC 120": D blank -> -13.6
$ws.Range("D30").Value = -13.6

# Row 31 "SC 132": D -13.7 -> blank; E blank -> -8.1
$ws.Range("D31").ClearContents()
$ws.Range("E31").Value = -8.1

# Row 32 "SC 193": D -14.7 -> blank
$ws.Range("D32").ClearContents()

# Row 33 "SC 232": E blank -> -10.7
$ws.Range("E33").Value = -10.7
